$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows 4-14: row-height tweaks ---
$ws.Rows.Item(4).RowHeight = 38.25
$ws.Rows.Item(5).RowHeight = 25.5
$ws.Rows.Item(6).RowHeight = 25.5
$ws.Rows.Item(7).RowHeight = 25.5
$ws.Rows.Item(8).RowHeight = 25.5
$ws.Rows.Item(9).RowHeight = 25.5
$ws.Rows.Item(10).RowHeight = 25.5
$ws.Rows.Item(11).RowHeight = 25.5
$ws.Rows.Item(12).RowHeight = 25.5
$ws.Rows.Item(13).RowHeight = 25.5
$ws.Rows.Item(14).RowHeight = 25.5

# --- Row 19: bold header row ---
$ws.Range("A19").Value = "Measuremen"
$ws.Range("B19").Value = "Token Prec"
$ws.Range("C19").Value = "Token Rec"
$ws.Range("D19").Value = "Token F1"
$ws.Range("E19").Value = "Span Pr"
$ws.Range("F19").Value = "Span Re"
$ws.Range("G19").Value = "Span F1"
$ws.Range("A19:G19").Font.Bold = $true

# --- Row 20 ---
$ws.Range("A20").Value = "CRF 10F"
$ws.Range("B20").Value = "vanilla w=4"
$ws.Range("C20").Value = "parserFE w=2"

# --- Row 21 ---
$ws.Range("A21").Value = "220 files"

# --- Rows 23-33: new results table (mirrors rows 4-14) ---
$ws.Range("A23").Value = "CRF + VanillaFE (window size 4)"
$ws.Range("B23").Value = 0.78876678876678796
$ws.Range("C23").Value = 0.39868339847767897
$ws.Range("D23").Value = 0.52965291063131903
$ws.Range("E23").Value = 0.68971848225214105
$ws.Range("F23").Value = 0.37170184696569902
$ws.Range("G23").Value = 0.483069009858551
$ws.Range("A23:G23").WrapText = $true
$ws.Rows.Item(23).RowHeight = 38.25

$ws.Range("A24").Value = "Test partition 1"
$ws.Range("B24").Value = 0.810276679841897
$ws.Range("C24").Value = 0.42797494780793299
$ws.Range("D24").Value = 0.56010928961748596
$ws.Range("E24").Value = 0.69047619047619002
$ws.Range("F24").Value = 0.4
$ws.Range("G24").Value = 0.50655021834061098
$ws.Range("A24:G24").WrapText = $true
$ws.Rows.Item(24).RowHeight = 25.5

$ws.Range("A25").Value = "Test partition 10"
$ws.Range("B25").Value = 0.76785714285714202
$ws.Range("C25").Value = 0.43877551020408101
$ws.Range("D25").Value = 0.55844155844155796
$ws.Range("E25").Value = 0.60526315789473595
$ws.Range("F25").Value = 0.40780141843971601
$ws.Range("G25").Value = 0.48728813559321998
$ws.Range("A25:G25").WrapText = $true
$ws.Rows.Item(25).RowHeight = 25.5

$ws.Range("A26").Value = "Test partition 2"
$ws.Range("B26").Value = 0.843537414965986
$ws.Range("C26").Value = 0.41891891891891803
$ws.Range("D26").Value = 0.55981941309255001
$ws.Range("E26").Value = 0.71
$ws.Range("F26").Value = 0.37765957446808501
$ws.Range("G26").Value = 0.49305555555555503
$ws.Range("A26:G26").WrapText = $true
$ws.Rows.Item(26).RowHeight = 25.5

$ws.Range("A27").Value = "Test partition 3"
$ws.Range("B27").Value = 0.75675675675675602
$ws.Range("C27").Value = 0.44444444444444398
$ws.Range("D27").Value = 0.56000000000000005
$ws.Range("E27").Value = 0.68531468531468498
$ws.Range("F27").Value = 0.4375
$ws.Range("G27").Value = 0.53405994550408697
$ws.Range("A27:G27").WrapText = $true
$ws.Rows.Item(27).RowHeight = 25.5

$ws.Range("A28").Value = "Test partition 4"
$ws.Range("B28").Value = 0.82033898305084696
$ws.Range("C28").Value = 0.456603773584905
$ws.Range("D28").Value = 0.586666666666666
$ws.Range("E28").Value = 0.75647668393782297
$ws.Range("F28").Value = 0.47096774193548302
$ws.Range("G28").Value = 0.58051689860834899
$ws.Range("A28:G28").WrapText = $true
$ws.Rows.Item(28).RowHeight = 25.5

$ws.Range("A29").Value = "Test partition 5"
$ws.Range("B29").Value = 0.77104377104377098
$ws.Range("C29").Value = 0.41335740072202098
$ws.Range("D29").Value = 0.53819036427732003
$ws.Range("E29").Value = 0.71052631578947301
$ws.Range("F29").Value = 0.38793103448275801
$ws.Range("G29").Value = 0.50185873605947895
$ws.Range("A29:G29").WrapText = $true
$ws.Rows.Item(29).RowHeight = 25.5

$ws.Range("A30").Value = "Test partition 6"
$ws.Range("B30").Value = 0.84523809523809501
$ws.Range("C30").Value = 0.41846758349705299
$ws.Range("D30").Value = 0.55978975032851497
$ws.Range("E30").Value = 0.77272727272727204
$ws.Range("F30").Value = 0.38746438746438699
$ws.Range("G30").Value = 0.51612903225806395
$ws.Range("A30:G30").WrapText = $true
$ws.Rows.Item(30).RowHeight = 25.5

$ws.Range("A31").Value = "Test partition 7"
$ws.Range("B31").Value = 0.81140350877192902
$ws.Range("C31").Value = 0.35238095238095202
$ws.Range("D31").Value = 0.49136786188579001
$ws.Range("E31").Value = 0.69444444444444398
$ws.Range("F31").Value = 0.32051282051281998
$ws.Range("G31").Value = 0.43859649122806998
$ws.Range("A31:G31").WrapText = $true
$ws.Rows.Item(31).RowHeight = 25.5

$ws.Range("A32").Value = "Test partition 8"
$ws.Range("B32").Value = 0.75
$ws.Range("C32").Value = 0.26844583987441101
$ws.Range("D32").Value = 0.39537572254335202
$ws.Range("E32").Value = 0.61783439490445802
$ws.Range("F32").Value = 0.234866828087167
$ws.Range("G32").Value = 0.34035087719298202
$ws.Range("A32:G32").WrapText = $true
$ws.Rows.Item(32).RowHeight = 25.5

$ws.Range("A33").Value = "Test partition 9"
$ws.Range("B33").Value = 0.72941176470588198
$ws.Range("C33").Value = 0.40172786177105801
$ws.Range("D33").Value = 0.51810584958217198
$ws.Range("E33").Value = 0.65317919075144504
$ws.Range("F33").Value = 0.35987261146496802
$ws.Range("G33").Value = 0.464065708418891
$ws.Range("A33:G33").WrapText = $true
$ws.Rows.Item(33).RowHeight = 25.5

# --- Final selection: D26 (matches target view state) ---
$ws.Range("D26").Select()
